# Swap the values of row 6 and row 7 for the columns that differ:
# A, B, E, F, G, H, Q, R (other columns already hold identical values
# in both rows, so they are unaffected either way).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell6 = $ws.Range($col + "6")
    $cell7 = $ws.Range($col + "7")

    $val6 = $cell6.Value2
    $val7 = $cell7.Value2

    $cell6.Value2 = $val7
    $cell7.Value2 = $val6
}
